$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 78). Bump each of those dates forward by one day
# (45189 -> 45190), matching the original cell's existing value so we
# don't disturb any row that might differ.
for ($r = 2; $r -le 78; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
